# "ekec + chem masters paths updated"
#
# The pivot-table sheet ("Sheet2") and its supporting pivot cache /
# pivot table definition are removed, leaving the module-listing sheet
# ("Sheet1") as the only worksheet in the workbook.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$pivotSheet = $wb.Worksheets.Item("Sheet2")

# Clear the pivot table itself first so Excel also drops the now-unused
# pivot cache / pivot table part from the package, then remove the sheet
# that hosted it.
$pivotTable = $pivotSheet.PivotTables().Item(1)
$pivotTable.TableRange2.Clear()
$pivotSheet.Delete()

$excel.DisplayAlerts = $true

# "Sheet1" is now the only sheet; make it active and leave the selection
# on H1 (one past the last data column).
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()
$ws.Range("H1").Select()
